$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04710838198661804
$ws.Range("C2").Value = 0.019473830237984657
$ws.Range("D2").Value = 0.015438690781593323
$ws.Range("E2").Value = 0.00775632169097662
$ws.Range("F2").Value = 0.000007074088443914661
$ws.Range("G2").Value = 0.002529376884922385
$ws.Range("H2").Value = 0.002292873105034232
$ws.Range("I2").Value = 1.2630733251571655
$ws.Range("J2").Value = 0.12570902705192566
$ws.Range("K2").Value = 1.4834043979644775
